# New weekly price report for "Agrícola del Norte S.A. de Arica - Camote".
# A new row of data is inserted above the current row 3, shifting the
# existing rows 3-5 down to rows 4-6 (all of their data is preserved as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, pushing rows 3..5 down to 4..6.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with this week's data.
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 45175
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 100114002
$ws.Range("G3").Value = "Camote"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 11500
$ws.Range("N3").Value = "$/malla 18 kilos"
$ws.Range("O3").Value = "Perú"
$ws.Range("P3").Value = 639
$ws.Range("Q3").Value = 18
$ws.Range("R3").Value = "Hortaliza"

# Match the date column's existing number format used by the rest of the
# column (D2, D4:D6 already carry it after the row insert).
$ws.Range("D3").NumberFormat = $ws.Range("D4").NumberFormat
